# Add a new "TextBox 26" shape (id=27) to the AGENDA slide (slide 3),
# listing the four agenda bullet points, matching the target OOXML diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The source deck's editing history churned through several add/undo
# cycles on this slide before the surviving textbox was created, so its
# shape id ended up at 27 ("TextBox 26") rather than the next free id
# (23). PowerPoint's per-slide shape-id counter keeps incrementing even
# across deletes, so replay that churn here to land on the same id.
for ($i = 0; $i -lt 4; $i++) {
    $dummy = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $dummy.Delete()
}

# AddTextbox's geometry arguments are in points; the target EMU values
# (712622, 1198102, 8221481, 2968761) divided by 12700 round-trips
# exactly back to those EMUs once PowerPoint re-serialises the shape.
$left   = 712622 / 12700
$top    = 1198102 / 12700
$width  = 8221481 / 12700
$height = 2968761 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 26"

$tb.Fill.Visible = 0

$tf2 = $tb.TextFrame2
$tf2.WordWrap = -1

$texts = @(
    "Introduction & Problem Statement: Overview of the project's objectives and challenges in photorealistic face generation.",
    "Solution Architecture: Explanation of the DCGAN architecture and its components.",
    "Training & Results: Description of the training procedure and presentation of results.",
    "Applications & Future Directions: Discussion of potential applications and future research directions."
)

$tr2 = $tf2.TextRange
$tr2.Text = $texts[0]
for ($i = 1; $i -lt $texts.Count; $i++) {
    $tr2.InsertAfter("`r" + $texts[$i]) | Out-Null
}

for ($i = 1; $i -le $texts.Count; $i++) {
    $para = $tr2.Paragraphs($i)

    $para.ParagraphFormat.SpaceWithin = 2

    $bullet = $para.ParagraphFormat.Bullet
    $bullet.Font.Name = "Arial"
    $bullet.Type = 1
    $bullet.Character = 8226

    $para.Font.Name = "Berlin Sans FB Demi"
    $para.Font.Size = 16
}

# spAutoFit recomputes the shape height from the rendered text; pin it
# back to the exact target extent afterwards. (Left/Top/Width are left
# untouched by AutoSize, and re-assigning them would lose precision on
# the float round-trip through the COM property getter/setter.)
$tf2.AutoSize = 1
$tb.Height = $height
